{"js": "// Replace the 25 \"NN\u00f7N=\" division prompts in the worksheet table with a\n// new set of prompts (one-to-one text substitution, each source string is\n// unique in the document so a plain search/replace is safe).\nconst replacements = [\n  [\"51\u00f76=\", \"63\u00f72=\"],\n  [\"17\u00f73=\", \"77\u00f77=\"],\n  [\"28\u00f78=\", \"82\u00f79=\"],\n  [\"76\u00f76=\", \"28\u00f79=\"],\n  [\"59\u00f79=\", \"84\u00f72=\"],\n  [\"97\u00f79=\", \"72\u00f77=\"],\n  [\"19\u00f73=\", \"34\u00f79=\"],\n  [\"81\u00f77=\", \"93\u00f79=\"],\n  [\"78\u00f78=\", \"13\u00f76=\"],\n  [\"78\u00f73=\", \"49\u00f79=\"],\n  [\"52\u00f78=\", \"58\u00f79=\"],\n  [\"81\u00f72=\", \"10\u00f77=\"],\n  [\"32\u00f76=\", \"68\u00f73=\"],\n  [\"33\u00f74=\", \"84\u00f74=\"],\n  [\"42\u00f76=\", \"80\u00f73=\"],\n  [\"14\u00f76=\", \"83\u00f77=\"],\n  [\"66\u00f74=\", \"31\u00f73=\"],\n  [\"83\u00f74=\", \"44\u00f74=\"],\n  [\"96\u00f79=\", \"37\u00f72=\"],\n  [\"43\u00f78=\", \"46\u00f75=\"],\n  [\"74\u00f75=\", \"49\u00f77=\"],\n  [\"80\u00f75=\", \"54\u00f79=\"],\n  [\"41\u00f75=\", \"63\u00f75=\"],\n  [\"88\u00f73=\", \"33\u00f76=\"],\n  [\"68\u00f74=\", \"58\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"NN\u00f7N=\" division prompts in the worksheet table with a\n# new set of prompts (one-to-one text substitution, each source string is\n# unique in the document so Find/Replace per pair is safe).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"51\u00f76=\", \"63\u00f72=\")\n  ,@(\"17\u00f73=\", \"77\u00f77=\")\n  ,@(\"28\u00f78=\", \"82\u00f79=\")\n  ,@(\"76\u00f76=\", \"28\u00f79=\")\n  ,@(\"59\u00f79=\", \"84\u00f72=\")\n  ,@(\"97\u00f79=\", \"72\u00f77=\")\n  ,@(\"19\u00f73=\", \"34\u00f79=\")\n  ,@(\"81\u00f77=\", \"93\u00f79=\")\n  ,@(\"78\u00f78=\", \"13\u00f76=\")\n  ,@(\"78\u00f73=\", \"49\u00f79=\")\n  ,@(\"52\u00f78=\", \"58\u00f79=\")\n  ,@(\"81\u00f72=\", \"10\u00f77=\")\n  ,@(\"32\u00f76=\", \"68\u00f73=\")\n  ,@(\"33\u00f74=\", \"84\u00f74=\")\n  ,@(\"42\u00f76=\", \"80\u00f73=\")\n  ,@(\"14\u00f76=\", \"83\u00f77=\")\n  ,@(\"66\u00f74=\", \"31\u00f73=\")\n  ,@(\"83\u00f74=\", \"44\u00f74=\")\n  ,@(\"96\u00f79=\", \"37\u00f72=\")\n  ,@(\"43\u00f78=\", \"46\u00f75=\")\n  ,@(\"74\u00f75=\", \"49\u00f77=\")\n  ,@(\"80\u00f75=\", \"54\u00f79=\")\n  ,@(\"41\u00f75=\", \"63\u00f75=\")\n  ,@(\"88\u00f73=\", \"33\u00f76=\")\n  ,@(\"68\u00f74=\", \"58\u00f77=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
